$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: Aptos -> PancakeSwap (rows re-ranked)
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.83"
$ws.Range("E31").Value = "  +7.64%  "

# Row 32: PancakeSwap -> Aptos (rows re-ranked)
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "6.67"
$ws.Range("E32").Value = "  +15.57%  "

# Row 46: Stellar -> Mantle (rows re-ranked)
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.595"
$ws.Range("E46").Value = "  +5.37%  "

# Row 47: Mantle -> Stellar (rows re-ranked)
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0961"
$ws.Range("E47").Value = "  +2.49%  "

# Remaining price / volume(1h) refreshes
$ws.Range("D2").Value = "63.389.95"
$ws.Range("E2").Value = "  +6.39%  "

$ws.Range("D3").Value = "2.437.96"
$ws.Range("E3").Value = "  +6.04%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "566.64"
$ws.Range("E5").Value = "  +4.64%  "

$ws.Range("D6").Value = "141.74"
$ws.Range("E6").Value = "  +10.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +4.18%  "

$ws.Range("D9").Value = "2.437.01"
$ws.Range("E9").Value = "  +6.08%  "

$ws.Range("E10").Value = "  +4.97%  "

$ws.Range("E11").Value = "  +4.13%  "

$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +7.53%  "

$ws.Range("D14").Value = "26.34"
$ws.Range("E14").Value = "  +14.15%  "

$ws.Range("D15").Value = "2.872.40"
$ws.Range("E15").Value = "  +6.09%  "

$ws.Range("D16").Value = "63.168.83"
$ws.Range("E16").Value = "  +6.30%  "

$ws.Range("D17").Value = "0.0000143"
$ws.Range("E17").Value = "  +9.22%  "

$ws.Range("D18").Value = "2.428.87"
$ws.Range("E18").Value = "  +5.99%  "

$ws.Range("D19").Value = "11.27"
$ws.Range("E19").Value = "  +8.57%  "

$ws.Range("D20").Value = "341.41"
$ws.Range("E20").Value = "  +10.31%  "

$ws.Range("D21").Value = "4.25"
$ws.Range("E21").Value = "  +5.77%  "

$ws.Range("D22").Value = "6.84"
$ws.Range("E22").Value = "  +5.37%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "65.35"
$ws.Range("E24").Value = "  +3.98%  "

$ws.Range("E25").Value = "  +3.85%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +14.87%  "

$ws.Range("D28").Value = "8.21"
$ws.Range("E28").Value = "  +6.74%  "

$ws.Range("D29").Value = "1.34"
$ws.Range("E29").Value = "  +13.33%  "

$ws.Range("E30").Value = "  +13.01%  "

$ws.Range("D33").Value = "174.13"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("E34").Value = "  +12.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.400"
$ws.Range("E35").Value = "  +6.24%  "

$ws.Range("D36").Value = "18.76"
$ws.Range("E36").Value = "  +6.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "373.30"
$ws.Range("E37").Value = "  +19.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.50"
$ws.Range("E38").Value = "  +13.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  +13.49%  "

$ws.Range("D42").Value = "39.95"
$ws.Range("E42").Value = "  +6.66%  "

$ws.Range("D43").Value = "148.58"
$ws.Range("E43").Value = "  +9.34%  "

$ws.Range("D44").Value = "3.71"
$ws.Range("E44").Value = "  +8.91%  "

$ws.Range("D45").Value = "20.88"
$ws.Range("E45").Value = "  +13.14%  "

$ws.Range("E48").Value = "  +7.33%  "

$ws.Range("E49").Value = "  +6.52%  "

$ws.Range("D50").Value = "17.92"
$ws.Range("E50").Value = "  +7.97%  "

$ws.Range("E51").Value = "  +17.31%  "
